$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A71").Value = 'S10'
$ws.Range("B71").Value = 'G01'
$ws.Range("C71").Value = 'Stock Universe & Groups: backend models and APIs'
$ws.Range("D71").Value = 'S10_G01_TB001'
$ws.Range("E71").Value = 'Define stocks, stock_groups, and stock_group_members models + migrations, and seed a minimal NSE equity universe.'
$ws.Range("F71").Value = 'See docs/stock_universe_group_backtests_prd.md for concept.'
$ws.Range("G71").Value = 'pending'

$ws.Range("A72").Value = 'S10'
$ws.Range("B72").Value = 'G01'
$ws.Range("C72").Value = 'Stock Universe & Groups: backend models and APIs'
$ws.Range("D72").Value = 'S10_G01_TB002'
$ws.Range("E72").Value = 'Expose REST APIs for stock universe CRUD and group membership management (list/create/update/delete, add/remove members).'
$ws.Range("F72").Value = 'Aligns Stocks page with backend; reuse existing FastAPI patterns.'
$ws.Range("G72").Value = 'pending'

$ws.Range("A73").Value = 'S10'
$ws.Range("B73").Value = 'G01'
$ws.Range("C73").Value = 'Stock Universe & Groups: backend models and APIs'
$ws.Range("D73").Value = 'S10_G01_TB003'
$ws.Range("E73").Value = 'Integrate universe/groups into BacktestService so group backtests can resolve symbols and validate coverage.'
$ws.Range("F73").Value = 'BacktestService can look up group symbols and warn when data is missing.'
$ws.Range("G73").Value = 'pending'

$ws.Range("A74").Value = 'S10'
$ws.Range("B74").Value = 'G02'
$ws.Range("C74").Value = 'Group backtests: portfolio simulator & capital allocation'
$ws.Range("D74").Value = 'S10_G02_TB001'
$ws.Range("E74").Value = 'Implement a portfolio simulator that consumes per-symbol candidate trades and enforces shared capital, max position size, per-trade risk, and broker constraints.'
$ws.Range("F74").Value = 'Initial policies: highestConfidenceSingle and allEligibleEqualWeight with default confidence=1.0.'
$ws.Range("G74").Value = 'pending'

$ws.Range("A75").Value = 'S10'
$ws.Range("B75").Value = 'G02'
$ws.Range("C75").Value = 'Group backtests: portfolio simulator & capital allocation'
$ws.Range("D75").Value = 'S10_G02_TB002'
$ws.Range("E75").Value = 'Extend BacktestService to run group backtests: load group members, run per-symbol strategy engines, feed candidates into portfolio simulator, and persist portfolio equity + trades.'
$ws.Range("F75").Value = 'Adds group_id/universe_mode to Backtest and keeps existing single-symbol path untouched.'
$ws.Range("G75").Value = 'pending'

$ws.Range("A76").Value = 'S10'
$ws.Range("B76").Value = 'G02'
$ws.Range("C76").Value = 'Group backtests: portfolio simulator & capital allocation'
$ws.Range("D76").Value = 'S10_G02_TB003'
$ws.Range("E76").Value = 'Compute portfolio-level realised/unrealised PnL and per-symbol summary metrics for group backtests.'
$ws.Range("F76").Value = 'Reuse existing PnL breakdown patterns and extend metrics_json schema.'
$ws.Range("G76").Value = 'pending'

$ws.Range("A77").Value = 'S10'
$ws.Range("B77").Value = 'G03'
$ws.Range("C77").Value = 'Stocks page UI: universe & groups management'
$ws.Range("D77").Value = 'S10_G03_TF001'
$ws.Range("E77").Value = 'Add Stocks page to sidebar with layout for Universe and Groups tabs.'
$ws.Range("F77").Value = 'Navigation only; no business logic change yet.'
$ws.Range("G77").Value = 'pending'

$ws.Range("A78").Value = 'S10'
$ws.Range("B78").Value = 'G03'
$ws.Range("C78").Value = 'Stocks page UI: universe & groups management'
$ws.Range("D78").Value = 'S10_G03_TF002'
$ws.Range("E78").Value = 'Implement Universe tab table and forms for adding/editing/deactivating stocks wired to /api/stocks APIs.'
$ws.Range("F78").Value = 'Follows existing Data/Strategies table patterns.'
$ws.Range("G78").Value = 'pending'

$ws.Range("A79").Value = 'S10'
$ws.Range("B79").Value = 'G03'
$ws.Range("C79").Value = 'Stocks page UI: universe & groups management'
$ws.Range("D79").Value = 'S10_G03_TF003'
$ws.Range("E79").Value = 'Implement Groups tab for creating/editing/deleting groups and managing group membership from the universe.'
$ws.Range("F79").Value = 'Provides source of truth for stock baskets such as trending_stocks.'
$ws.Range("G79").Value = 'pending'

$ws.Range("A80").Value = 'S10'
$ws.Range("B80").Value = 'G04'
$ws.Range("C80").Value = 'Backtests UI: group runs & portfolio reporting'
$ws.Range("D80").Value = 'S10_G04_TF001'
$ws.Range("E80").Value = 'Extend Run Backtest form with target selector (single stock vs stock group) and group dropdown, updating payloads for group backtests.'
$ws.Range("F80").Value = 'Group runs share risk/cost settings and initial capital across all symbols.'
$ws.Range("G80").Value = 'pending'

$ws.Range("A81").Value = 'S10'
$ws.Range("B81").Value = 'G04'
$ws.Range("C81").Value = 'Backtests UI: group runs & portfolio reporting'
$ws.Range("D81").Value = 'S10_G04_TF002'
$ws.Range("E81").Value = 'Update Backtest Details to display group context, portfolio PnL breakdown (realised/unrealised), and per-symbol summary table for group backtests.'
$ws.Range("F81").Value = 'Keeps existing single-symbol layout while adding portfolio view when group_id is present.'
$ws.Range("G81").Value = 'pending'

$ws.Range("A82").Value = 'S10'
$ws.Range("B82").Value = 'G04'
$ws.Range("C82").Value = 'Backtests UI: group runs & portfolio reporting'
$ws.Range("D82").Value = 'S10_G04_TF003'
$ws.Range("E82").Value = 'Ensure trades table and CSV export behave as a portfolio trade ledger (symbol-aware) for group runs.'
$ws.Range("F82").Value = 'Reuses existing trades export format with symbol column and new group context.'
$ws.Range("G82").Value = 'pending'
